$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E19").Value = "SMD 1/4W 470"
$ws.Range("G19").Value = 0.03
$ws.Range("H19").Value = "https://www.reichelt.de/SMD-1206-von-0-bis-910-Ohm/SMD-1-4W-470/3/index.html?ACTION=3&GROUPID=7973&ARTICLE=18337&OFFSET=75&"

$ws.Range("C20").Value = 1
$ws.Range("F20").ClearContents()
$ws.Range("G20").Value = 0.02

$ws.Range("E33").Value = "from malectrics.eu"

$ws.Range("E35").Select()
